$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = 729270
$ws.Range("D13").Value = 734126
$ws.Range("D14").Value = 5184815

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D14"), 0, 2)
$ws.Sort.SetRange($ws.Range("A2:E14"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

$ws.Range("G4").Formula = "=E9/E4"
$ws.Range("G4").NumberFormat = "#,##0"
$ws.Range("G5").Formula = "=E14/E4"

$ws.Range("A8").Select()
